$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are purely numeric-looking strings (e.g. "587.18").
# The source data stores them as text, so force a Text number format on
# those specific cells first to prevent Excel from auto-converting the
# assigned string into a numeric value.
$textCells = 'D5','D6','D13','D19','D23','D24','D26','D27','D34','D35','D39','D43','D45','D46'
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = '@'
}

$ws.Range('D2').Value = '63.139.07'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.569.52'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('D5').Value = '587.18'
$ws.Range('E5').Value = '  +3.26%  '
$ws.Range('D6').Value = '148.38'
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('E8').Value = '  +1.65%  '
$ws.Range('E9').Value = '  +2.70%  '
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = '27.64'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').Value = '3.032.97'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').Value = '63.022.57'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('D17').Value = '2.575.17'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '343.62'
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('E20').Value = '  +2.82%  '
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '5.53'
$ws.Range('E23').Value = '  -3.72%  '
$ws.Range('D24').Value = '66.60'
$ws.Range('E24').Value = '  +1.99%  '
$ws.Range('D25').Value = '2.670.43'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.170'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').Value = '1.63'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').Value = '  +11.38%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('E32').Value = '  +7.01%  '
$ws.Range('D33').Value = '0.0₃0827'
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('D34').Value = '463.48'
$ws.Range('E34').Value = '  +12.31%  '
$ws.Range('D35').Value = '176.59'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('E36').Value = '  +2.69%  '
$ws.Range('E37').Value = '  +1.18%  '
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('D39').Value = '4.63'
$ws.Range('E39').Value = '  +4.97%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '151.34'
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('E44').Value = '  +1.23%  '
$ws.Range('D45').Value = '21.06'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('D46').Value = '0.0551'
$ws.Range('E46').Value = '  +5.18%  '
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('E51').Value = '  +0.45%  '
